# Weekly fruit/vegetable price update: a new "Feria Lagunitas de Puerto Montt"
# Melón entry (date 2023-01-05 / serial 44931) is inserted as two rows
# (Calameño + Tuna, both "Primera") right after the existing row 325,
# shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 326:327, pushing the old 326.. rows down to 328..
$ws.Rows("326:327").Insert()

# New row 326 - Calameño
$ws.Range("A326").Value = 4
$ws.Range("B326").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C326").Value = "Los Lagos"
$ws.Range("D326").Value = 44931
$ws.Range("E326").Value = 10
$ws.Range("F326").Value = 100112027
$ws.Range("G326").Value = "Melón"
$ws.Range("H326").Value = "Calameño"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 6000
$ws.Range("K326").Value = 2000
$ws.Range("L326").Value = 2000
$ws.Range("M326").Value = 2000
$ws.Range("N326").Value = "$/unidad"
$ws.Range("O326").Value = "Región del Maule"
$ws.Range("P326").Value = 2000
$ws.Range("Q326").Value = 1
$ws.Range("R326").Value = "Hortaliza"

# New row 327 - Tuna
$ws.Range("A327").Value = 4
$ws.Range("B327").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C327").Value = "Los Lagos"
$ws.Range("D327").Value = 44931
$ws.Range("E327").Value = 10
$ws.Range("F327").Value = 100112027
$ws.Range("G327").Value = "Melón"
$ws.Range("H327").Value = "Tuna"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 6000
$ws.Range("K327").Value = 2000
$ws.Range("L327").Value = 2000
$ws.Range("M327").Value = 2000
$ws.Range("N327").Value = "$/unidad"
$ws.Range("O327").Value = "Región del Maule"
$ws.Range("P327").Value = 2000
$ws.Range("Q327").Value = 1
$ws.Range("R327").Value = "Hortaliza"
